$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 40.73489498692954
$ws.Range("B2").Value = 0.92990332116438268
$ws.Range("C2").Value = 0.32236534800031658
$ws.Range("D2").Value = -0.17708866619884475

$ws.Range("A3").Value = -4.4358889630384315
$ws.Range("B3").Value = 0.3342830573024545
$ws.Range("C3").Value = -0.5399207721196182
$ws.Range("D3").Value = 0.77248973937152032

$ws.Range("A4").Value = -81.554208193969302
$ws.Range("B4").Value = 0.15341007427146044
$ws.Range("C4").Value = -0.77753851495754855
$ws.Range("D4").Value = -0.60983539325758629

$ws.Range("E4").Select()
